$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. AgreementInfo: remove the "DOT" test-case row (was row 4) so only
#    the ALT and ROW rows remain.
# ------------------------------------------------------------------
$agreementInfo = $wb.Worksheets.Item("AgreementInfo")
$agreementInfo.Rows.Item(4).Delete()
$agreementInfo.Range("B7").Select()

# ------------------------------------------------------------------
# 2. Add a new "PayeeInfo" worksheet right after "AgreementInfo" and
#    populate it with the AddPayeeInformationALT test-case row.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$payeeInfo = $wb.Worksheets.Add($null, $lastSheet)
$payeeInfo.Name = "PayeeInfo"

$payeeInfo.Range("A1").Value = "String"
$payeeInfo.Range("B1").Value = "String"
$payeeInfo.Range("C1").Value = "String"

$payeeInfo.Range("A2").Value = "Title"
$payeeInfo.Range("A3").Value = "AddPayeeInformationALT"
$payeeInfo.Range("B2").Value = "LandownerName"
$payeeInfo.Range("B3").Value = "Long Form Renewable"
$payeeInfo.Range("C3").Value = "Parcel #: test1, Grantor Name: , County PID:"
$payeeInfo.Range("C2").Value = "AvailableTract"

# Match the header/title formatting already used on the other sheets.
$agreementInfo.Range("A2").Copy()
$payeeInfo.Range("A2").PasteSpecial(-4122)
$agreementInfo.Range("B2").Copy()
$payeeInfo.Range("B2").PasteSpecial(-4122)
$agreementInfo.Range("C2").Copy()
$payeeInfo.Range("C2").PasteSpecial(-4122)
$agreementInfo.Range("B3").Copy()
$payeeInfo.Range("B3").PasteSpecial(-4122)

$payeeInfo.Columns.Item(1).ColumnWidth = 27.95
$payeeInfo.Columns.Item(2).ColumnWidth = 19.95
$payeeInfo.Columns.Item(3).ColumnWidth = 39.15

$payeeInfo.Range("C6").Select()
